$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Update the value for "coloured_arteries" levels/range from [1, 2, 3] to [0, 1, 2, 3]
$ws.Range("C13").Value = "[0, 1, 2, 3]"

# Update the last active selection to match the edited cell
$ws.Range("C13").Select()
